$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.880.21"
$ws.Range("E2").Value = "  +1.68%  "

# Row 3
$ws.Range("D3").Value = "2.499.14"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "

# Row 8
$ws.Range("E8").Value = "  -0.30%  "

# Row 9
$ws.Range("D9").Value = "2.497.75"
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("E10").Value = "  +5.35%  "

# Row 11
$ws.Range("E11").Value = "  -1.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "

# Row 14
$ws.Range("D14").Value = "2.951.83"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.93%  "

# Row 16
$ws.Range("D16").Value = "68.731.58"
$ws.Range("E16").Value = "  +1.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000172"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").Value = "2.495.00"
$ws.Range("E18").Value = "  -0.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "360.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "

# Row 23
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.99%  "

# Row 27
$ws.Range("B27").Value = "SuiNetwork"
$ws.Range("C27").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.25%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.625.99"
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "502.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.92%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0868"
$ws.Range("E31").Value = "  -4.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.87%  "

# Row 33
$ws.Range("E33").Value = "  -0.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.90%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.54%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.316"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.69%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.81%  "

# Row 49
$ws.Range("E49").Value = "  -1.56%  "

# Row 50
$ws.Range("E50").Value = "  -2.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.574"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.14%  "
